$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.809.10"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "3.506.38"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.38%  "
$ws.Range("D9").Value = "3.512.88"
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "4.115.83"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("D17").Value = "65.807.91"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").Value = "3.492.67"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.553"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000125"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.41%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("D40").Value = "3.093.71"
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0775"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0324"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "315.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.55%  "
